$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.131.19"
$ws.Range("E2").Value = "  -0.92%  "

$ws.Range("D3").Value = "1.769.44"
$ws.Range("E3").Value = "  -1.07%  "

$ws.Range("D4").Value = "'1.008"

$ws.Range("D5").Value = "'333.49"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("D7").Value = "'0.3762"
$ws.Range("E7").Value = "  -0.93%  "

$ws.Range("D8").Value = "'0.3390"
$ws.Range("E8").Value = "  -3.03%  "

$ws.Range("D9").Value = "'47.91"
$ws.Range("E9").Value = "  -3.46%  "

$ws.Range("D10").Value = "'1.177"
$ws.Range("E10").Value = "  -3.71%  "

$ws.Range("D11").Value = "'0.07331"
$ws.Range("E11").Value = "  -4.13%  "

$ws.Range("D12").Value = "'1.005"
$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("D13").Value = "'21.30"
$ws.Range("E13").Value = "  -1.44%  "

$ws.Range("D14").Value = "'6.350"
$ws.Range("E14").Value = "  -4.11%  "

$ws.Range("D15").Value = "1.771.32"
$ws.Range("E15").Value = "  -0.92%  "

$ws.Range("D16").Value = "'6.955"
$ws.Range("E16").Value = "  -3.64%  "

$ws.Range("D17").Value = "'0.00001076"
$ws.Range("E17").Value = "  -3.80%  "

$ws.Range("D18").Value = "'0.06642"
$ws.Range("E18").Value = "  -2.23%  "

$ws.Range("D19").Value = "'83.40"
$ws.Range("E19").Value = "  -2.43%  "

$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Value = "'6.519"
$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").Value = "'17.06"
$ws.Range("E22").Value = "  -3.33%  "

$ws.Range("D23").Value = "27.104.68"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("D24").Value = "'12.28"
$ws.Range("E24").Value = "  -6.80%  "

$ws.Range("D25").Value = "'2.419"
$ws.Range("E25").Value = "  -1.83%  "

$ws.Range("D26").Value = "'1.490"
$ws.Range("E26").Value = "  -3.57%  "

$ws.Range("D27").Value = "'2.488"
$ws.Range("E27").Value = "  -2.94%  "

$ws.Range("D28").Value = "'20.88"
$ws.Range("E28").Value = "  +2.62%  "

$ws.Range("D29").Value = "'150.37"
$ws.Range("E29").Value = "  -2.15%  "

$ws.Range("D30").Value = "1.974.03"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").Value = "'132.14"
$ws.Range("E31").Value = "  -2.62%  "

$ws.Range("D32").Value = "'4.059"
$ws.Range("E32").Value = "  -2.91%  "

$ws.Range("D33").Value = "'5.907"
$ws.Range("E33").Value = "  -7.76%  "

$ws.Range("D34").Value = "'0.08543"
$ws.Range("E34").Value = "  -2.62%  "

$ws.Range("D35").Value = "'12.85"
$ws.Range("E35").Value = "  -5.20%  "

$ws.Range("D36").Value = "'1.652"
$ws.Range("E36").Value = "  -4.37%  "

$ws.Range("D37").Value = "'5.352"
$ws.Range("E37").Value = "  -5.26%  "

$ws.Range("D38").Value = "'0.6726"
$ws.Range("E38").Value = "  -1.31%  "

$ws.Range("D39").Value = "'0.06276"
$ws.Range("E39").Value = "  -3.98%  "

$ws.Range("D40").Value = "'0.02319"
$ws.Range("E40").Value = "  -4.49%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.661"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.2158"
$ws.Range("E42").Value = "  -4.95%  "

$ws.Range("D43").Value = "'1.228"
$ws.Range("E43").Value = "  -1.27%  "

$ws.Range("D44").Value = "'14.30"
$ws.Range("E44").Value = "  -3.64%  "

$ws.Range("D45").Value = "'1.006"
$ws.Range("E45").Value = "  +0.63%  "

$ws.Range("D46").Value = "'0.6254"
$ws.Range("E46").Value = "  -3.18%  "

$ws.Range("D47").Value = "'3.817"
$ws.Range("E47").Value = "  -3.99%  "

$ws.Range("D48").Value = "'2.094"
$ws.Range("E48").Value = "  -3.15%  "

$ws.Range("D49").Value = "'128.38"
$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").Value = "'0.07114"
$ws.Range("E50").Value = "  -3.46%  "

$ws.Range("D51").Value = "'78.11"
$ws.Range("E51").Value = "  -3.04%  "
